$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 410
$ws.Range("I2").Value = 150
$ws.Range("J2").Value = 583.3333
$ws.Range("K2").Value = 150
$ws.Range("L2").Value = 583.3333
$ws.Range("M2").Value = -37
$ws.Range("N2").Value = -809.3333

# Row 6
$ws.Range("H6").Value = 1933.2727
$ws.Range("I6").Value = 88.5
$ws.Range("J6").Value = 2625.0625
$ws.Range("K6").Value = 265.5
$ws.Range("L6").Value = 7875.1875
$ws.Range("M6").Value = -153.5
$ws.Range("N6").Value = -8099.1875

# Row 29
$ws.Range("H29").Value = 3114.7144
$ws.Range("I29").Value = 400.6
$ws.Range("J29").Value = 9900
$ws.Range("K29").Value = 1201.8
$ws.Range("L29").Value = 29700
$ws.Range("M29").Value = -920.8000000000002
$ws.Range("N29").Value = -30262

# Row 38
$ws.Range("H38").Value = 524.25
$ws.Range("I38").Value = 524.25
$ws.Range("K38").Value = 1572.75
$ws.Range("M38").Value = -1200.75

# Row 43
$ws.Range("H43").Value = 1318
$ws.Range("I43").Value = 1080
$ws.Range("J43").Value = 1377.5
$ws.Range("K43").Value = 1080
$ws.Range("L43").Value = 1377.5
$ws.Range("M43").Value = -1011
$ws.Range("N43").Value = -1515.5

# Row 62
$ws.Range("H62").Value = 2957.1428
$ws.Range("I62").Value = 2933.3333
$ws.Range("K62").Value = 2933.3333
$ws.Range("M62").Value = -2309.3333

# Row 65
$ws.Range("H65").Value = 2957.1428
$ws.Range("I65").Value = 2933.3333
$ws.Range("K65").Value = 14666.6665
$ws.Range("M65").Value = -11546.6665

# Row 137
$ws.Range("H137").Value = 1736.2963
$ws.Range("I137").Value = 1550.875
$ws.Range("K137").Value = 4652.625
$ws.Range("M137").Value = -2102.625

# Row 141
$ws.Range("H141").Value = 780076.6
$ws.Range("I141").Value = 967303.5600000001
$ws.Range("K141").Value = 2901910.68
$ws.Range("M141").Value = -2896730.68

$ws = $wb.Worksheets.Item("ARM")
# Row 11
$ws.Range("H11").Value = 6690001.5
$ws.Range("J11").Value = 70004
$ws.Range("L11").Value = 70004
$ws.Range("N11").Value = -70292

# Row 45
$ws.Range("H45").Value = 1566.6923
$ws.Range("I45").Value = 1179
$ws.Range("K45").Value = 1179
$ws.Range("M45").Value = -802

# Row 74
$ws.Range("H74").Value = 1551.5
$ws.Range("I74").Value = 1410.3846
$ws.Range("J74").Value = 1718.2727
$ws.Range("K74").Value = 1410.3846
$ws.Range("L74").Value = 1718.2727
$ws.Range("M74").Value = -536.3846000000001
$ws.Range("N74").Value = -3466.2727

# Row 77
$ws.Range("H77").Value = 1551.5
$ws.Range("I77").Value = 1410.3846
$ws.Range("J77").Value = 1718.2727
$ws.Range("K77").Value = 7051.923000000001
$ws.Range("L77").Value = 8591.363499999999
$ws.Range("M77").Value = -2683.923000000001
$ws.Range("N77").Value = -17327.3635

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 128207.25
$ws.Range("I86").Value = 2750.8333
$ws.Range("J86").Value = 203481.1
$ws.Range("K86").Value = 2750.8333
$ws.Range("L86").Value = 203481.1
$ws.Range("M86").Value = -1627.8333
$ws.Range("N86").Value = -205727.1

# Row 89
$ws.Range("H89").Value = 128207.25
$ws.Range("I89").Value = 2750.8333
$ws.Range("J89").Value = 203481.1
$ws.Range("K89").Value = 13754.1665
$ws.Range("L89").Value = 1017405.5
$ws.Range("M89").Value = -8138.166499999999
$ws.Range("N89").Value = -1028637.5

# Row 107
$ws.Range("H107").Value = 5333.6665
$ws.Range("I107").Value = 5333.6665
$ws.Range("K107").Value = 5333.6665
$ws.Range("M107").Value = -3413.6665

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 823.4
$ws.Range("I16").Value = 776
$ws.Range("K16").Value = 776
$ws.Range("M16").Value = -489

# Row 31
$ws.Range("H31").Value = 1738.3334
$ws.Range("I31").Value = 1505.9231
$ws.Range("J31").Value = 2342.6
$ws.Range("K31").Value = 1505.9231
$ws.Range("L31").Value = 2342.6
$ws.Range("M31").Value = -1210.9231
$ws.Range("N31").Value = -2932.6

# Row 34
$ws.Range("H34").Value = 1738.3334
$ws.Range("I34").Value = 1505.9231
$ws.Range("J34").Value = 2342.6
$ws.Range("K34").Value = 1505.9231
$ws.Range("L34").Value = 2342.6
$ws.Range("M34").Value = -1303.9231
$ws.Range("N34").Value = -2746.6

# Row 99
$ws.Range("H99").Value = 2841.4
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502

# Row 113
$ws.Range("H113").Value = 823.4
$ws.Range("I113").Value = 776
$ws.Range("K113").Value = 776
$ws.Range("M113").Value = 1394

# Row 126
$ws.Range("H126").Value = 2841.4
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530

$ws = $wb.Worksheets.Item("CUL")
# Row 97
$ws.Range("H97").Value = 971.3333
$ws.Range("J97").Value = 971.3333
$ws.Range("L97").Value = 2913.9999
$ws.Range("N97").Value = -3905.9999

# Row 114
$ws.Range("H114").Value = 2152.6667
$ws.Range("J114").Value = 2377.6
$ws.Range("L114").Value = 7132.799999999999
$ws.Range("N114").Value = -13640.8

# Row 131
$ws.Range("H131").Value = 9205.905000000001
$ws.Range("J131").Value = 9890.966
$ws.Range("L131").Value = 29672.898
$ws.Range("N131").Value = -39752.898

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 30366.666
$ws.Range("I70").Value = 43500
$ws.Range("K70").Value = 43500
$ws.Range("M70").Value = -43230

# Row 73
$ws.Range("H73").Value = 30366.666
$ws.Range("I73").Value = 43500
$ws.Range("K73").Value = 43500
$ws.Range("M73").Value = -42564

# Row 113
$ws.Range("H113").Value = 1474.5
$ws.Range("I113").Value = 1449
$ws.Range("K113").Value = 1449
$ws.Range("M113").Value = 721

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 3900.8333
$ws.Range("I16").Value = 6587.4
$ws.Range("J16").Value = 1981.8572
$ws.Range("K16").Value = 6587.4
$ws.Range("L16").Value = 1981.8572
$ws.Range("M16").Value = -6417.4
$ws.Range("N16").Value = -2321.8572

# Row 22
$ws.Range("H22").Value = 1930
$ws.Range("I22").Value = 1803.3334
$ws.Range("J22").Value = 2500
$ws.Range("K22").Value = 1803.3334
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = -1508.3334
$ws.Range("N22").Value = -3090

# Row 27
$ws.Range("H27").Value = 1930
$ws.Range("I27").Value = 1803.3334
$ws.Range("J27").Value = 2500
$ws.Range("K27").Value = 1803.3334
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = -1696.3334
$ws.Range("N27").Value = -2714

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1248.025
$ws.Range("I132").Value = 988.96155
$ws.Range("K132").Value = 2966.88465
$ws.Range("M132").Value = -436.88465
